$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.40288149
$ws.Range("H2").Value = 373.98178118
$ws.Range("M2").Value = 0.5502426597222885
$ws.Range("N2").Value = 55.91611202929847

$ws.Range("G3").Value = 6.78317677
$ws.Range("H3").Value = 610.8127203800001
$ws.Range("M3").Value = 0.6792581408243431
$ws.Range("N3").Value = 114.4994539192471

$ws.Range("G4").Value = 2.53737894
$ws.Range("H4").Value = 71.83620431
$ws.Range("M4").Value = 0.3919814959048554
$ws.Range("N4").Value = 21.28889590035775

$ws.Range("G5").Value = 2.79238079
$ws.Range("H5").Value = 134.33426993
$ws.Range("M5").Value = 0.3614862467633796
$ws.Range("N5").Value = 33.534404364962

$ws.Range("G6").Value = 0.8456443100000001
$ws.Range("H6").Value = 12.93354327
$ws.Range("M6").Value = 0.2371911856015121
$ws.Range("N6").Value = 6.212888359056365

$ws.Range("G7").Value = 1.00181138
$ws.Range("H7").Value = 25.66063576
$ws.Range("M7").Value = 0.1950750737072204
$ws.Range("N7").Value = 8.607412685903103

$ws.Range("G8").Value = 0.41613435
$ws.Range("H8").Value = 4.070704539999999
$ws.Range("M8").Value = 0.1450668131519251
$ws.Range("N8").Value = 2.221703299735914

$ws.Range("G9").Value = 0.48189639
$ws.Range("H9").Value = 8.538840759999999
$ws.Range("M9").Value = 0.1188358353462049
$ws.Range("N9").Value = 3.785142394448388

$ws.Range("G10").Value = 0.2117835
$ws.Range("H10").Value = 1.53689199
$ws.Range("M10").Value = 0.09005140704320205
$ws.Range("N10").Value = 0.970451232026567

$ws.Range("G11").Value = 0.25371109
$ws.Range("H11").Value = 3.546879979999999
$ws.Range("M11").Value = 0.07577188831599005
$ws.Range("N11").Value = 2.081288556471711

$ws.Range("G12").Value = 0.12865917
$ws.Range("H12").Value = 0.7916509599999999
$ws.Range("M12").Value = 0.06378650906225743
$ws.Range("N12").Value = 0.5961162616030753

$ws.Range("G13").Value = 0.14540667
$ws.Range("H13").Value = 1.63893001
$ws.Range("M13").Value = 0.04765479951408701
$ws.Range("N13").Value = 0.9543881149747753
